# Adds a new row (row 8) of checkout/shipping test data to the active sheet,
# matching the row layout already used by rows 2-7 (A:H = Address, Email,
# FirstName, LastName, Phone, City, State, ZipCode).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "25 f2 Address (2)"
$ws.Range("B8").Value = " "
$ws.Range("C8").Value = "25 f2 First (1)"
$ws.Range("D8").Value = "25 F2 Last (4)"
# Numeric-looking values must stay text, so prefix with an apostrophe.
$ws.Range("E8").Value = "'2502834949444"
$ws.Range("F8").Value = "25 f2 City ( 3)"
$ws.Range("G8").Value = "'252"
$ws.Range("H8").Value = "'250205"
